# Scheduled runner update: refresh Universalis market price snapshots
# and derived Leve profit columns (H:N) across all job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 2475.5557
$ws.Range("I51").Value = 1175
$ws.Range("J51").Value = 3516
$ws.Range("K51").Value = 1175
$ws.Range("L51").Value = 3516
$ws.Range("M51").Value = -691
$ws.Range("N51").Value = -4484
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H70").Value = 5810.091
$ws.Range("I70").Value = 25780.5
$ws.Range("J70").Value = 1372.2222
$ws.Range("K70").Value = 77341.5
$ws.Range("L70").Value = 4116.6666
$ws.Range("M70").Value = -77071.5
$ws.Range("N70").Value = -4656.6666
$ws.Range("H73").Value = 5810.091
$ws.Range("I73").Value = 25780.5
$ws.Range("J73").Value = 1372.2222
$ws.Range("K73").Value = 77341.5
$ws.Range("L73").Value = 4116.6666
$ws.Range("M73").Value = -76405.5
$ws.Range("N73").Value = -5988.6666
$ws.Range("H107").Value = 7326.353
$ws.Range("I107").Value = 11080.091
$ws.Range("J107").Value = 444.5
$ws.Range("K107").Value = 11080.091
$ws.Range("L107").Value = 444.5
$ws.Range("M107").Value = -9160.091
$ws.Range("N107").Value = -4284.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 21277.158
$ws.Range("I32").Value = 21477.74
$ws.Range("J32").Value = 17666.666
$ws.Range("K32").Value = 21477.74
$ws.Range("L32").Value = 17666.666
$ws.Range("M32").Value = -21190.74
$ws.Range("N32").Value = -18240.666
$ws.Range("H37").Value = 7732.8335
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 7732.8335
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 7732.8335
$ws.Range("M37").ClearContents()
$ws.Range("N37").Value = -8278.833500000001
$ws.Range("H55").Value = 31991.8
$ws.Range("J55").Value = 39239.75
$ws.Range("L55").Value = 39239.75
$ws.Range("N55").Value = -39869.75
$ws.Range("H63").Value = 29414762
$ws.Range("I63").Value = 41669330
$ws.Range("J63").Value = 3799.6
$ws.Range("K63").Value = 41669330
$ws.Range("L63").Value = 3799.6
$ws.Range("M63").Value = -41668644
$ws.Range("N63").Value = -5171.6
$ws.Range("H66").Value = 29414762
$ws.Range("I66").Value = 41669330
$ws.Range("J66").Value = 3799.6
$ws.Range("K66").Value = 208346650
$ws.Range("L66").Value = 18998
$ws.Range("M66").Value = -208343218
$ws.Range("N66").Value = -25862
$ws.Range("H80").Value = 29107.5
$ws.Range("J80").Value = 29107.5
$ws.Range("L80").Value = 29107.5
$ws.Range("N80").Value = -31103.5
$ws.Range("H83").Value = 29107.5
$ws.Range("J83").Value = 29107.5
$ws.Range("L83").Value = 87322.5
$ws.Range("N83").Value = -97306.5
$ws.Range("H102").Value = 1543.375
$ws.Range("I102").Value = 1308.1666
$ws.Range("J102").Value = 2249
$ws.Range("K102").Value = 1308.1666
$ws.Range("L102").Value = 2249
$ws.Range("M102").Value = 313.8334
$ws.Range("N102").Value = -5493
$ws.Range("H132").Value = 4316.8887
$ws.Range("I132").Value = 1971.3529
$ws.Range("J132").Value = 11566.728
$ws.Range("K132").Value = 5914.0587
$ws.Range("L132").Value = 34700.18399999999
$ws.Range("M132").Value = -3384.0587
$ws.Range("N132").Value = -39760.18399999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 60000
$ws.Range("J35").Value = 60000
$ws.Range("L35").Value = 60000
$ws.Range("N35").Value = -60620
$ws.Range("H82").Value = 26999.889
$ws.Range("J82").Value = 32434.6
$ws.Range("L82").Value = 32434.6
$ws.Range("N82").Value = -33200.6
$ws.Range("H85").Value = 26999.889
$ws.Range("J85").Value = 32434.6
$ws.Range("L85").Value = 32434.6
$ws.Range("N85").Value = -35086.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4250.342
$ws.Range("I31").Value = 1437.1666
$ws.Range("J31").Value = 8606.226000000001
$ws.Range("K31").Value = 1437.1666
$ws.Range("L31").Value = 8606.226000000001
$ws.Range("M31").Value = -1142.1666
$ws.Range("N31").Value = -9196.226000000001
$ws.Range("H34").Value = 4250.342
$ws.Range("I34").Value = 1437.1666
$ws.Range("J34").Value = 8606.226000000001
$ws.Range("K34").Value = 1437.1666
$ws.Range("L34").Value = 8606.226000000001
$ws.Range("M34").Value = -1235.1666
$ws.Range("N34").Value = -9010.226000000001
$ws.Range("H105").Value = 1782.5454
$ws.Range("I105").Value = 1810.8
$ws.Range("J105").Value = 1500
$ws.Range("K105").Value = 1810.8
$ws.Range("L105").Value = 1500
$ws.Range("M105").Value = -63.79999999999995
$ws.Range("N105").Value = -4994
$ws.Range("H132").Value = 4316.8887
$ws.Range("I132").Value = 1971.3529
$ws.Range("J132").Value = 11566.728
$ws.Range("K132").Value = 5914.0587
$ws.Range("L132").Value = 34700.18399999999
$ws.Range("M132").Value = -3384.0587
$ws.Range("N132").Value = -39760.18399999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H125").Value = 4489
$ws.Range("J125").Value = 4773.8335
$ws.Range("L125").Value = 14321.5005
$ws.Range("N125").Value = -24161.5005
$ws.Range("H132").Value = 1799.3684
$ws.Range("I132").Value = 1904
$ws.Range("J132").Value = 1793.5555
$ws.Range("K132").Value = 17136
$ws.Range("L132").Value = 16141.9995
$ws.Range("M132").Value = -14606
$ws.Range("N132").Value = -21201.9995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 772.6445
$ws.Range("I97").Value = 697.8333
$ws.Range("J97").Value = 922.26666
$ws.Range("K97").Value = 697.8333
$ws.Range("L97").Value = 922.26666
$ws.Range("M97").Value = -201.8333
$ws.Range("N97").Value = -1914.26666
$ws.Range("H126").Value = 3861.389
$ws.Range("I126").Value = 2672.5715
$ws.Range("K126").Value = 8017.7145
$ws.Range("M126").Value = -5547.7145

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3434.8965
$ws.Range("I40").Value = 4311.647
$ws.Range("J40").Value = 2192.8333
$ws.Range("K40").Value = 4311.647
$ws.Range("L40").Value = 2192.8333
$ws.Range("M40").Value = -4175.647
$ws.Range("N40").Value = -2464.8333
$ws.Range("H46").Value = 961.25
$ws.Range("I46").Value = 470
$ws.Range("J46").Value = 1256
$ws.Range("K46").Value = 470
$ws.Range("L46").Value = 1256
$ws.Range("M46").Value = -282
$ws.Range("N46").Value = -1632
$ws.Range("H68").Value = 2503
$ws.Range("I68").Value = 2078.5715
$ws.Range("J68").Value = 3351.8572
$ws.Range("K68").Value = 2078.5715
$ws.Range("L68").Value = 3351.8572
$ws.Range("M68").Value = -1329.5715
$ws.Range("N68").Value = -4849.8572
$ws.Range("H71").Value = 2503
$ws.Range("I71").Value = 2078.5715
$ws.Range("J71").Value = 3351.8572
$ws.Range("K71").Value = 10392.8575
$ws.Range("L71").Value = 16759.286
$ws.Range("M71").Value = -6648.8575
$ws.Range("N71").Value = -24247.286
$ws.Range("H93").Value = 1613.65
$ws.Range("I93").Value = 1662.0714
$ws.Range("J93").Value = 1500.6666
$ws.Range("K93").Value = 1662.0714
$ws.Range("L93").Value = 1500.6666
$ws.Range("M93").Value = -414.0714
$ws.Range("N93").Value = -3996.6666
$ws.Range("H122").Value = 6134.852
$ws.Range("I122").Value = 5016.933
$ws.Range("J122").Value = 7532.25
$ws.Range("K122").Value = 15050.799
$ws.Range("L122").Value = 22596.75
$ws.Range("M122").Value = -12600.799
$ws.Range("N122").Value = -27496.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1337.5883
$ws.Range("I136").Value = 945.6585
$ws.Range("J136").Value = 2944.5
$ws.Range("K136").Value = 2836.9755
$ws.Range("L136").Value = 8833.5
$ws.Range("M136").Value = -286.9755
$ws.Range("N136").Value = -13933.5
